$d = $word.ActiveDocument

# Locate the "Docker Compose" heading run and append an "e" right after it,
# with the same run-level formatting (bold, underline, size 28), turning the
# heading text into "Docker Composee".

$rng = $d.Content
$found = $rng.Find.Execute("Docker Compose", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

$rng.Collapse(0)  # wdCollapseEnd
$rng.InsertAfter("e")
$rng.Font.Bold = $true
$rng.Font.Underline = 1  # wdUnderlineSingle
$rng.Font.Size = 14
